$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (data for target cluster combos that no longer exist)
$ws.Rows("8:10").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf12"
$ws.Range("C2").Value = "Cd163"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.115481333333334
$ws.Range("H2").Value = 15.346444
$ws.Range("I2").Value = 0.1917470154127355
$ws.Range("J2").Value = 0.1917470154127354
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01985533333333333
$ws.Range("N2").Value = 0.059566
$ws.Range("O2").Value = 0.02825030780058695
$ws.Range("P2").Value = 0.02825030780058696
$ws.Range("Q2").Value = 0.1015695870337778
$ws.Range("R2").Value = 0.9141262833040001
$ws.Range("S2").Value = 0.005416912205253667
$ws.Range("T2").Value = 0.005416912205253668

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf12"
$ws.Range("C3").Value = "Cd163"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.115481333333334
$ws.Range("H3").Value = 15.346444
$ws.Range("I3").Value = 0.1917470154127355
$ws.Range("J3").Value = 0.1917470154127354
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6829806666666668
$ws.Range("N3").Value = 2.048942
$ws.Range("O3").Value = 0.9717496921994131
$ws.Range("P3").Value = 0.9717496921994131
$ws.Range("Q3").Value = 3.49377485136089
$ws.Range("R3").Value = 31.44397366224801
$ws.Range("S3").Value = 0.1863301032074818
$ws.Range("T3").Value = 0.1863301032074818

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf12"
$ws.Range("C4").Value = "Cd163"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.971037
$ws.Range("H4").Value = 35.913111
$ws.Range("I4").Value = 0.4487184033275903
$ws.Range("J4").Value = 0.4487184033275903
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01985533333333333
$ws.Range("N4").Value = 0.059566
$ws.Range("O4").Value = 0.02825030780058695
$ws.Range("P4").Value = 0.02825030780058696
$ws.Range("Q4").Value = 0.2376889299806667
$ws.Range("R4").Value = 2.139200369826
$ws.Range("S4").Value = 0.01267643300979235
$ws.Range("T4").Value = 0.01267643300979235

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf12"
$ws.Range("C5").Value = "Cd163"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.971037
$ws.Range("H5").Value = 35.913111
$ws.Range("I5").Value = 0.4487184033275903
$ws.Range("J5").Value = 0.4487184033275903
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6829806666666668
$ws.Range("N5").Value = 2.048942
$ws.Range("O5").Value = 0.9717496921994131
$ws.Range("P5").Value = 0.9717496921994131
$ws.Range("Q5").Value = 8.175986830951336
$ws.Range("R5").Value = 73.58388147856201
$ws.Range("S5").Value = 0.436041970317798
$ws.Range("T5").Value = 0.436041970317798

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Tnfsf12"
$ws.Range("C6").Value = "Cd163"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.591765666666667
$ws.Range("H6").Value = 28.775297
$ws.Range("I6").Value = 0.3595345812596742
$ws.Range("J6").Value = 0.3595345812596742
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01985533333333333
$ws.Range("N6").Value = 0.059566
$ws.Range("O6").Value = 0.02825030780058695
$ws.Range("P6").Value = 0.02825030780058696
$ws.Range("Q6").Value = 0.1904477045668889
$ws.Range("R6").Value = 1.714029341102
$ws.Range("S6").Value = 0.01015696258554094
$ws.Range("T6").Value = 0.01015696258554094

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Tnfsf12"
$ws.Range("C7").Value = "Cd163"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.591765666666667
$ws.Range("H7").Value = 28.775297
$ws.Range("I7").Value = 0.3595345812596742
$ws.Range("J7").Value = 0.3595345812596742
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6829806666666668
$ws.Range("N7").Value = 2.048942
$ws.Range("O7").Value = 0.9717496921994131
$ws.Range("P7").Value = 0.9717496921994131
$ws.Range("Q7").Value = 6.550990509530446
$ws.Range("R7").Value = 58.95891458577401
$ws.Range("S7").Value = 0.3493776186741333
$ws.Range("T7").Value = 0.3493776186741333
